# LOB1256.docx content shuffle
# Applies the reordering of the "Objetivos / Docente(s) / Programa resumido /
# Programa / Avaliacao / Bibliografia" sections described by the target diff.
#
# Strategy: capture every text block that needs to move BEFORE any mutation
# (paragraph offsets earlier in the document stay stable as long as we only
# edit paragraphs at-or-after the current one), then rewrite the paragraphs
# from the END of the document back to the START so that earlier offsets
# we already captured/compute remain valid.

$d = $word.ActiveDocument
$vt = [char]11   # Word's in-paragraph line-break char (<w:br/>) inside Range.Text

# ---------------------------------------------------------------------------
# 1) Capture source text blocks, before touching anything.
# ---------------------------------------------------------------------------

# Para 6: "Fornecer os fundamentos ..." (Objetivos, PT)
$p6 = $d.Paragraphs.Item(6)
$txt_fornecer = $d.Range($p6.Range.Start, $p6.Range.End - 1).Text

# Para 7: "Provide the fundamentals ..." (Objetivos, EN - italic)
$p7 = $d.Paragraphs.Item(7)
$txt_provide = $d.Range($p7.Range.Start, $p7.Range.End - 1).Text

# Para 9: "9146830 - Danubia ..." <br/> "5464150 - Mariana ..." (Docentes)
$p9 = $d.Paragraphs.Item(9)
$p9full = $d.Range($p9.Range.Start, $p9.Range.End - 1).Text
$p9parts = $p9full.Split($vt)
$txt_danubia = $p9parts[0]
$txt_mariana = $p9parts[1]

# Para 11: PT short program summary
$p11 = $d.Paragraphs.Item(11)
$txt_intro_pt = $d.Range($p11.Range.Start, $p11.Range.End - 1).Text

# Para 12: EN short program summary (italic)
$p12 = $d.Paragraphs.Item(12)
$txt_intro_en = $d.Range($p12.Range.Start, $p12.Range.End - 1).Text

# Para 14: PT full "Programa" text
$p14 = $d.Paragraphs.Item(14)
$txt_programa_pt = $d.Range($p14.Range.Start, $p14.Range.End - 1).Text

# Para 17: Avaliacao bullet paragraph - 3 bold labels + 3 value runs.
# Use Find to split on the (unchanged, unique-within-paragraph) bold labels
# rather than walking Runs directly.
$p17 = $d.Paragraphs.Item(17)
$pStart = $p17.Range.Start
$pEndExclMark = $p17.Range.End - 1

$fr = $d.Range($pStart, $pEndExclMark)
$fr.Find.ClearFormatting()
$fr.Find.Text = "M" + [char]233 + "todo: "
$fr.Find.Execute() | Out-Null
$metodoLabelEnd = $fr.End

$fr2 = $d.Range($metodoLabelEnd, $pEndExclMark)
$fr2.Find.ClearFormatting()
$fr2.Find.Text = "Crit" + [char]233 + "rio: "
$fr2.Find.Execute() | Out-Null
$criterioLabelStart = $fr2.Start
$criterioLabelEnd = $fr2.End

$fr3 = $d.Range($criterioLabelEnd, $pEndExclMark)
$fr3.Find.ClearFormatting()
$fr3.Find.Text = "Norma de recupera" + [char]231 + [char]227 + "o: "
$fr3.Find.Execute() | Out-Null
$normaLabelStart = $fr3.Start
$normaLabelEnd = $fr3.End

$txt_metodo_value = $d.Range($metodoLabelEnd, $criterioLabelStart).Text
$txt_criterio_value = $d.Range($criterioLabelEnd, $normaLabelStart).Text
$txt_norma_value = $d.Range($normaLabelEnd, $pEndExclMark).Text

# Para 19: full Bibliografia text block
$p19 = $d.Paragraphs.Item(19)
$txt_bibliografia = $d.Range($p19.Range.Start, $p19.Range.End - 1).Text

# ---------------------------------------------------------------------------
# 2) Rewrite paragraphs, from the last one back to the first, so offsets
#    captured/used for earlier paragraphs stay valid.
# ---------------------------------------------------------------------------

# Para 19 -> "5464150 - Mariana Consiglio Kasemodel"
$p19b = $d.Paragraphs.Item(19)
$d.Range($p19b.Range.Start, $p19b.Range.End - 1).Text = $txt_mariana

# Para 17 -> rebuild the 6-run Avaliacao bullet with swapped content:
#   Metodo:       (unchanged label)
#   <norma value> + trailing break
#   Criterio:     (unchanged label)
#   <bibliografia text> + trailing break
#   Norma de recuperacao: (unchanged label, now after the bibliography)
#   <Danubia line>
$p17b = $d.Paragraphs.Item(17)
$pStart = $p17b.Range.Start
$pEndExclMark = $p17b.Range.End - 1
$newPara17 = "M" + [char]233 + "todo: " + $vt + $txt_norma_value + $vt + "Crit" + [char]233 + "rio: " + $vt + $txt_bibliografia + $vt + "Norma de recupera" + [char]231 + [char]227 + "o: " + $vt + $txt_danubia
$d.Range($pStart, $pEndExclMark).Text = $newPara17

# Re-apply bold to the three labels (the plain .Text rewrite collapsed
# everything into a single run, so the bold formatting on the labels needs
# to be restored).
$p17c = $d.Paragraphs.Item(17)
$searchBase = $d.Range($p17c.Range.Start, $p17c.Range.End - 1)

$fb1 = $d.Range($p17c.Range.Start, $p17c.Range.End - 1)
$fb1.Find.ClearFormatting()
$fb1.Find.Text = "M" + [char]233 + "todo: "
$fb1.Find.Execute() | Out-Null
$fb1.Bold = 1

$fb2 = $d.Range($p17c.Range.Start, $p17c.Range.End - 1)
$fb2.Find.ClearFormatting()
$fb2.Find.Text = "Crit" + [char]233 + "rio: "
$fb2.Find.Execute() | Out-Null
$fb2.Bold = 1

$fb3 = $d.Range($p17c.Range.Start, $p17c.Range.End - 1)
$fb3.Find.ClearFormatting()
$fb3.Find.Text = "Norma de recupera" + [char]231 + [char]227 + "o: "
$fb3.Find.Execute() | Out-Null
$fb3.Bold = 1

# Para 14 -> "Media ponderada de exercicios e provas." (drop trailing break)
$p14b = $d.Paragraphs.Item(14)
$criterioValueClean = $txt_criterio_value.TrimEnd($vt)
$d.Range($p14b.Range.Start, $p14b.Range.End - 1).Text = $criterioValueClean

# Para 12 -> "Provide the fundamentals ..." (stays italic automatically)
$p12b = $d.Paragraphs.Item(12)
$d.Range($p12b.Range.Start, $p12b.Range.End - 1).Text = $txt_provide

# Para 11 -> "Aulas expositivas ..." (drop trailing break)
$p11b = $d.Paragraphs.Item(11)
$metodoValueClean = $txt_metodo_value.TrimEnd($vt)
$d.Range($p11b.Range.Start, $p11b.Range.End - 1).Text = $metodoValueClean

# Para 9 -> "Fornecer os fundamentos ..." <br/> <full PT Programa text>
$p9b = $d.Paragraphs.Item(9)
$newPara9 = $txt_fornecer + $vt + $txt_programa_pt
$d.Range($p9b.Range.Start, $p9b.Range.End - 1).Text = $newPara9

# Para 7 -> EN short program summary (stays italic automatically)
$p7b = $d.Paragraphs.Item(7)
$d.Range($p7b.Range.Start, $p7b.Range.End - 1).Text = $txt_intro_en

# Para 6 -> PT short program summary
$p6b = $d.Paragraphs.Item(6)
$d.Range($p6b.Range.Start, $p6b.Range.End - 1).Text = $txt_intro_pt

Write-Host "Done."
